$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[53.58048282556054, 75.09308573772826]"
$ws.Range("N2").Value = [double]"1.110223024625157e-15"
$ws.Range("O2").Value = [double]"1.110223024625157e-15"
$ws.Range("U2").Value = "[41.60189324682019, 55.722562874216784]"
$ws.Range("M3").Value = "[52.32905360881331, 76.65329993459859]"
$ws.Range("N3").Value = [double]"6.372680161348399e-14"
$ws.Range("O3").Value = [double]"6.372680161348399e-14"
$ws.Range("U3").Value = "[41.80825094652131, 54.64061128974295]"
$ws.Range("M4").Value = "[50.994569563484326, 75.42002950769412]"
$ws.Range("N4").Value = [double]"1.394440118929197e-13"
$ws.Range("O4").Value = [double]"1.394440118929197e-13"
$ws.Range("Q4").Value = "[2.232763547596349, 2.6352899336418885]"
$ws.Range("U4").Value = "[42.95017422401657, 55.613554417099586]"
$ws.Range("Y4").Value = [double]"13.02822822822827"
$ws.Range("Z4").Value = [double]"14.46582582582587"
$ws.Range("M5").Value = "[51.557497779586654, 73.07524382345032]"
$ws.Range("N5").Value = [double]"3.33066907387547e-15"
$ws.Range("O5").Value = [double]"3.33066907387547e-15"
$ws.Range("U5").Value = "[45.54916183479121, 58.30691561260768]"
$ws.Range("M6").Value = "[54.051948793386785, 75.26966523332455]"
$ws.Range("N6").Value = [double]"4.440892098500626e-16"
$ws.Range("O6").Value = [double]"4.440892098500626e-16"
$ws.Range("Q6").Value = "[-3.207632138800389, -2.8554215510105427]"
$ws.Range("U6").Value = "[43.296498742852684, 56.46536100901785]"
$ws.Range("Y6").Value = [double]"11.27955955955976"
$ws.Range("Z6").Value = [double]"12.6708708708711"
$ws.Range("M7").Value = "[53.79558405227728, 75.0147385461499]"
$ws.Range("N7").Value = [double]"6.661338147750939e-16"
$ws.Range("O7").Value = [double]"6.661338147750939e-16"
$ws.Range("U7").Value = "[43.38875679044257, 56.59522317179092]"
$ws.Range("M8").Value = "[52.414958214833305, 75.75366518530494]"
$ws.Range("N8").Value = [double]"2.020605904817785e-14"
$ws.Range("O8").Value = [double]"2.020605904817785e-14"
$ws.Range("U8").Value = "[41.448996085265655, 54.717343248135684]"
$ws.Range("M9").Value = "[51.24177447198315, 75.99357490229144]"
$ws.Range("N9").Value = [double]"1.731947918415244e-13"
$ws.Range("O9").Value = [double]"1.731947918415244e-13"
$ws.Range("U9").Value = "[43.792685424885406, 57.04554346350113]"
$ws.Range("M10").Value = "[50.30214559764306, 75.98520746141794]"
$ws.Range("N10").Value = [double]"7.036593530074242e-13"
$ws.Range("O10").Value = [double]"7.036593530074242e-13"
$ws.Range("Q10").Value = "[2.144710900648888, 2.5472372866944273]"
$ws.Range("U10").Value = "[44.26469901206232, 57.437146642255186]"
$ws.Range("M11").Value = "[50.426938291187746, 75.61373754846252]"
$ws.Range("N11").Value = [double]"4.061195824078823e-13"
$ws.Range("O11").Value = [double]"4.061195824078823e-13"
$ws.Range("U11").Value = "[44.5743326628551, 57.627937199223865]"
$ws.Range("M12").Value = "[51.778023852495295, 74.50728709348462]"
$ws.Range("N12").Value = [double]"1.354472090042691e-14"
$ws.Range("O12").Value = [double]"1.354472090042691e-14"
$ws.Range("U12").Value = "[45.34635770212569, 57.95474898118475]"
$ws.Range("M13").Value = "[53.65485833027368, 73.43007431160146]"
$ws.Range("N13").Value = [double]"2.220446049250313e-16"
$ws.Range("O13").Value = [double]"2.220446049250313e-16"
$ws.Range("U13").Value = "[43.89577799219937, 56.104140483456604]"
$ws.Range("M14").Value = "[53.99447608976868, 72.4103140127709]"
$ws.Range("U14").Value = "[44.007522716351716, 56.57792603107703]"

Write-Host "Done applying changes"
